$d = $word.ActiveDocument

# --- Step 1: register numbering definitions (abstractNum + num) matching a fresh "1." decimal list ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchRange = $scratchPara.Range
$scratchRange.Text = "scratch"
$scratchRange.ListFormat.ApplyNumberDefault()
$scratchPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara2.Range.Delete()

# --- Step 2: find the paragraph containing 10/15/2024 and insert a clean paragraph right after it ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "10/15/2024") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate 10/15/2024 paragraph"
}
$target.Range.InsertParagraphAfter()
$targetIndex = $target.Range.Information(3)
$newPara = $target.Next()
$r = $newPara.Range
$r.Collapse(1)
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Next Steps for Implementation:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Implement Teacher Finalization Feature</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Develop AI-Driven Scheduling Suggestions</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Google Calendar Integration</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Cancellation System with Optional Messaging</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Missed Lesson Tracking and Post-Lesson Notifications</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Post-Confirmation Error Handling</w:t></w:r></w:p>'
$r.InsertXML($xmlFrag)

Write-Host "count: $($d.Paragraphs.Count)"
for ($i = 99; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    Write-Host "$i : len=$($t.Length) [$t]"
}
